$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Sheet1 previously carried 43 leftover "index-only" rows (45-87) beyond the
# real 14-column data block (which ends at row 44). Remove them so the
# sheet's used range / dimension shrinks back to A1:N44.
$ws1.Range("A45:A87").EntireRow.Delete()

# The active tab moves from Sheet3 back to Sheet1 (bug fix - Sheet3 was left
# selected), with a new scroll position/selection sitting on Sheet1.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F61").Select()
